# Auto-generated Excel COM-interop edit script
# Applies cell-level numeric updates to columns H-N across multiple sheets
# (currentAveragePrice, currentAveragePriceNQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 371  # H2: 340.2 -> 371
$ws.Cells.Item(2, 9).Value = 182.9  # I2: 162.91667 -> 182.9
$ws.Cells.Item(2, 11).Value = 182.9  # K2: 162.91667 -> 182.9
$ws.Cells.Item(2, 13).Value = -69.90000000000001  # M2: -49.91667000000001 -> -69.90000000000001
$ws.Cells.Item(5, 8).Value = 175.375  # H5: 158 -> 175.375
$ws.Cells.Item(5, 9).Value = 186.14285  # I5: 165.25 -> 186.14285
$ws.Cells.Item(5, 11).Value = 186.14285  # K5: 165.25 -> 186.14285
$ws.Cells.Item(5, 13).Value = -71.14285000000001  # M5: -50.25 -> -71.14285000000001
$ws.Cells.Item(6, 8).Value = 16.25  # H6: 18.333334 -> 16.25
$ws.Cells.Item(6, 9).Value = 16.25  # I6: 18.333334 -> 16.25
$ws.Cells.Item(6, 11).Value = 48.75  # K6: 55.000002 -> 48.75
$ws.Cells.Item(6, 13).Value = 63.25  # M6: 56.999998 -> 63.25
$ws.Cells.Item(9, 8).Value = 53.4  # H9: 54.81818 -> 53.4
$ws.Cells.Item(9, 9).Value = 60.375  # I9: 61.333332 -> 60.375
$ws.Cells.Item(9, 11).Value = 60.375  # K9: 61.333332 -> 60.375
$ws.Cells.Item(9, 13).Value = 108.625  # M9: 107.666668 -> 108.625
$ws.Cells.Item(12, 8).Value = 455  # H12: 468.8 -> 455
$ws.Cells.Item(12, 9).Value = 343.75  # I12: 361 -> 343.75
$ws.Cells.Item(12, 11).Value = 343.75  # K12: 361 -> 343.75
$ws.Cells.Item(12, 13).Value = -173.75  # M12: -191 -> -173.75
$ws.Cells.Item(18, 8).Value = 1937.5  # H18: 1558.3334 -> 1937.5
$ws.Cells.Item(18, 10).Value = 3000  # J18: 1900 -> 3000
$ws.Cells.Item(18, 12).Value = 3000  # L18: 1900 -> 3000
$ws.Cells.Item(18, 14).Value = -3568  # N18: -2468 -> -3568
$ws.Cells.Item(29, 8).Value = 43.5  # H29: 44 -> 43.5
$ws.Cells.Item(29, 9).Value = 43.5  # I29: 44 -> 43.5
$ws.Cells.Item(29, 11).Value = 130.5  # K29: 132 -> 130.5
$ws.Cells.Item(29, 13).Value = 150.5  # M29: 149 -> 150.5
$ws.Cells.Item(33, 8).Value = 269.9  # H33: 250.36363 -> 269.9
$ws.Cells.Item(33, 9).Value = 267.22223  # I33: 246 -> 267.22223
$ws.Cells.Item(33, 11).Value = 267.22223  # K33: 246 -> 267.22223
$ws.Cells.Item(33, 13).Value = -38.22223000000002  # M33: -17 -> -38.22223000000002
$ws.Cells.Item(41, 8).Value = 950.2  # H41: 824.5 -> 950.2
$ws.Cells.Item(41, 9).Value = 100.5  # I41: 132.33333 -> 100.5
$ws.Cells.Item(41, 11).Value = 100.5  # K41: 132.33333 -> 100.5
$ws.Cells.Item(41, 13).Value = 339.5  # M41: 307.66667 -> 339.5
$ws.Cells.Item(53, 8).Value = 1048.9375  # H53: 988.05884 -> 1048.9375
$ws.Cells.Item(53, 9).Value = 697  # I53: 628.7 -> 697
$ws.Cells.Item(53, 11).Value = 697  # K53: 628.7 -> 697
$ws.Cells.Item(53, 13).Value = -60  # M53: 8.299999999999955 -> -60
$ws.Cells.Item(55, 8).Value = 62.625  # H55: 61.625 -> 62.625
$ws.Cells.Item(55, 9).Value = 91.14286  # I55: 81.5 -> 91.14286
$ws.Cells.Item(55, 10).Value = 40.444443  # J55: 41.75 -> 40.444443
$ws.Cells.Item(55, 11).Value = 91.14286  # K55: 81.5 -> 91.14286
$ws.Cells.Item(55, 12).Value = 40.444443  # L55: 41.75 -> 40.444443
$ws.Cells.Item(55, 13).Value = 122.85714  # M55: 132.5 -> 122.85714
$ws.Cells.Item(55, 14).Value = -468.444443  # N55: -469.75 -> -468.444443
$ws.Cells.Item(98, 8).Value = 2153.889  # H98: 2154 -> 2153.889
$ws.Cells.Item(98, 9).Value = 2153.889  # I98: 2154 -> 2153.889
$ws.Cells.Item(98, 11).Value = 2153.889  # K98: 2154 -> 2153.889
$ws.Cells.Item(98, 13).Value = -655.8890000000001  # M98: -656 -> -655.8890000000001
$ws.Cells.Item(107, 8).Value = 228.27272  # H107: 223.82608 -> 228.27272
$ws.Cells.Item(107, 9).Value = 139.5625  # I107: 138.76471 -> 139.5625
$ws.Cells.Item(107, 11).Value = 139.5625  # K107: 138.76471 -> 139.5625
$ws.Cells.Item(107, 13).Value = 1780.4375  # M107: 1781.23529 -> 1780.4375
$ws.Cells.Item(112, 8).Value = 2229  # H112: 2050.25 -> 2229
$ws.Cells.Item(112, 10).Value = 2660.8  # J112: 2350.5 -> 2660.8
$ws.Cells.Item(112, 12).Value = 7982.400000000001  # L112: 7051.5 -> 7982.400000000001
$ws.Cells.Item(112, 14).Value = -10198.4  # N112: -9267.5 -> -10198.4
$ws.Cells.Item(122, 8).Value = 2153.889  # H122: 2154 -> 2153.889
$ws.Cells.Item(122, 9).Value = 2153.889  # I122: 2154 -> 2153.889
$ws.Cells.Item(122, 11).Value = 6461.667  # K122: 6462 -> 6461.667
$ws.Cells.Item(122, 13).Value = -4011.667  # M122: -4012 -> -4011.667
$ws.Cells.Item(138, 8).Value = 4198.3  # H138: 3882.2 -> 4198.3
$ws.Cells.Item(138, 10).Value = 4499.25  # J138: 4018.7693 -> 4499.25
$ws.Cells.Item(138, 12).Value = 13497.75  # L138: 12056.3079 -> 13497.75
$ws.Cells.Item(138, 14).Value = -23777.75  # N138: -22336.3079 -> -23777.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(8, 8).Value = 2175  # H8: 0 -> 2175
$ws.Cells.Item(8, 9).Value = 500  # I8: 0 -> 500
$ws.Cells.Item(8, 10).Value = 3850  # J8: 0 -> 3850
$ws.Cells.Item(8, 11).Value = 500  # K8: 0 -> 500
$ws.Cells.Item(8, 12).Value = 3850  # L8: 0 -> 3850
$ws.Cells.Item(8, 13).Value = -356  # M8: None -> -356
$ws.Cells.Item(8, 14).Value = -4138  # N8: None -> -4138
$ws.Cells.Item(13, 8).Value = 1607.8334  # H13: 2999 -> 1607.8334
$ws.Cells.Item(13, 9).Value = 216.66667  # I13: 0 -> 216.66667
$ws.Cells.Item(13, 11).Value = 216.66667  # K13: 0 -> 216.66667
$ws.Cells.Item(13, 13).Value = -72.66667000000001  # M13: None -> -72.66667000000001
$ws.Cells.Item(137, 8).Value = 0  # H137: 65000 -> 0
$ws.Cells.Item(137, 10).Value = 0  # J137: 65000 -> 0
$ws.Cells.Item(137, 12).Value = 0  # L137: 65000 -> 0
$ws.Cells.Item(137, 14).ClearContents()  # N137: -75200 -> (removed)
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 9599.571  # H20: 9000.556 -> 9599.571
$ws.Cells.Item(20, 9).Value = 9829.333000000001  # I20: 8639.200000000001 -> 9829.333000000001
$ws.Cells.Item(20, 10).Value = 9427.25  # J20: 9452.25 -> 9427.25
$ws.Cells.Item(20, 11).Value = 9829.333000000001  # K20: 8639.200000000001 -> 9829.333000000001
$ws.Cells.Item(20, 12).Value = 9427.25  # L20: 9452.25 -> 9427.25
$ws.Cells.Item(20, 13).Value = -9582.333000000001  # M20: -8392.200000000001 -> -9582.333000000001
$ws.Cells.Item(20, 14).Value = -9921.25  # N20: -9946.25 -> -9921.25
$ws.Cells.Item(59, 8).Value = 94999  # H59: 0 -> 94999
$ws.Cells.Item(59, 10).Value = 94999  # J59: 0 -> 94999
$ws.Cells.Item(59, 12).Value = 94999  # L59: 0 -> 94999
$ws.Cells.Item(59, 14).Value = -96693  # N59: None -> -96693
$ws.Cells.Item(80, 8).Value = 213.88235  # H80: 219.4375 -> 213.88235
$ws.Cells.Item(80, 10).Value = 220.25  # J80: 228.90909 -> 220.25
$ws.Cells.Item(80, 12).Value = 220.25  # L80: 228.90909 -> 220.25
$ws.Cells.Item(80, 14).Value = -2216.25  # N80: -2224.90909 -> -2216.25
$ws.Cells.Item(83, 8).Value = 213.88235  # H83: 219.4375 -> 213.88235
$ws.Cells.Item(83, 10).Value = 220.25  # J83: 228.90909 -> 220.25
$ws.Cells.Item(83, 12).Value = 1101.25  # L83: 1144.54545 -> 1101.25
$ws.Cells.Item(83, 14).Value = -11085.25  # N83: -11128.54545 -> -11085.25
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 9840.682000000001  # H31: 9760.652 -> 9840.682000000001
$ws.Cells.Item(31, 9).Value = 9562.5  # I31: 9388.888999999999 -> 9562.5
$ws.Cells.Item(31, 11).Value = 9562.5  # K31: 9388.888999999999 -> 9562.5
$ws.Cells.Item(31, 13).Value = -9267.5  # M31: -9093.888999999999 -> -9267.5
$ws.Cells.Item(34, 8).Value = 9840.682000000001  # H34: 9760.652 -> 9840.682000000001
$ws.Cells.Item(34, 9).Value = 9562.5  # I34: 9388.888999999999 -> 9562.5
$ws.Cells.Item(34, 11).Value = 9562.5  # K34: 9388.888999999999 -> 9562.5
$ws.Cells.Item(34, 13).Value = -9360.5  # M34: -9186.888999999999 -> -9360.5
$ws.Cells.Item(107, 8).Value = 542.8  # H107: 561.619 -> 542.8
$ws.Cells.Item(107, 9).Value = 503.1111  # I107: 526 -> 503.1111
$ws.Cells.Item(107, 11).Value = 503.1111  # K107: 526 -> 503.1111
$ws.Cells.Item(107, 13).Value = 1416.8889  # M107: 1394 -> 1416.8889
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(19, 8).Value = 99  # H19: 100 -> 99
$ws.Cells.Item(19, 9).Value = 99  # I19: 100 -> 99
$ws.Cells.Item(19, 11).Value = 297  # K19: 300 -> 297
$ws.Cells.Item(19, 13).Value = -123  # M19: -126 -> -123
$ws.Cells.Item(48, 8).Value = 0  # H48: 387.75 -> 0
$ws.Cells.Item(48, 9).Value = 0  # I48: 183.33333 -> 0
$ws.Cells.Item(48, 10).Value = 0  # J48: 1001 -> 0
$ws.Cells.Item(48, 11).Value = 0  # K48: 549.99999 -> 0
$ws.Cells.Item(48, 12).Value = 0  # L48: 3003 -> 0
$ws.Cells.Item(48, 13).ClearContents()  # M48: -299.99999 -> (removed)
$ws.Cells.Item(48, 14).ClearContents()  # N48: -3503 -> (removed)
$ws.Cells.Item(69, 8).Value = 0  # H69: 2500 -> 0
$ws.Cells.Item(69, 10).Value = 0  # J69: 2500 -> 0
$ws.Cells.Item(69, 12).Value = 0  # L69: 7500 -> 0
$ws.Cells.Item(69, 14).ClearContents()  # N69: -9122 -> (removed)
$ws.Cells.Item(72, 8).Value = 0  # H72: 2500 -> 0
$ws.Cells.Item(72, 10).Value = 0  # J72: 2500 -> 0
$ws.Cells.Item(72, 12).Value = 0  # L72: 22500 -> 0
$ws.Cells.Item(72, 14).ClearContents()  # N72: -30612 -> (removed)
$ws.Cells.Item(80, 8).Value = 4788.636  # H80: 4797.826 -> 4788.636
$ws.Cells.Item(81, 8).Value = 3013  # H81: 3049.5 -> 3013
$ws.Cells.Item(81, 10).Value = 3013  # J81: 3049.5 -> 3013
$ws.Cells.Item(81, 12).Value = 9039  # L81: 9148.5 -> 9039
$ws.Cells.Item(81, 14).Value = -11285  # N81: -11394.5 -> -11285
$ws.Cells.Item(83, 8).Value = 4788.636  # H83: 4797.826 -> 4788.636
$ws.Cells.Item(84, 8).Value = 3013  # H84: 3049.5 -> 3013
$ws.Cells.Item(84, 10).Value = 3013  # J84: 3049.5 -> 3013
$ws.Cells.Item(84, 12).Value = 27117  # L84: 27445.5 -> 27117
$ws.Cells.Item(84, 14).Value = -38349  # N84: -38677.5 -> -38349
$ws.Cells.Item(88, 8).Value = 14999.5  # H88: 15000 -> 14999.5
$ws.Cells.Item(88, 10).Value = 14999.5  # J88: 15000 -> 14999.5
$ws.Cells.Item(88, 12).Value = 44998.5  # L88: 45000 -> 44998.5
$ws.Cells.Item(88, 14).Value = -45854.5  # N88: -45856 -> -45854.5
$ws.Cells.Item(91, 8).Value = 14999.5  # H91: 15000 -> 14999.5
$ws.Cells.Item(91, 10).Value = 14999.5  # J91: 15000 -> 14999.5
$ws.Cells.Item(91, 12).Value = 44998.5  # L91: 45000 -> 44998.5
$ws.Cells.Item(91, 14).Value = -47962.5  # N91: -47964 -> -47962.5
$ws.Cells.Item(98, 8).Value = 371.30768  # H98: 387.58334 -> 371.30768
$ws.Cells.Item(98, 9).Value = 223.57143  # I98: 224.14285 -> 223.57143
$ws.Cells.Item(98, 10).Value = 543.6667  # J98: 616.4 -> 543.6667
$ws.Cells.Item(98, 11).Value = 670.71429  # K98: 672.4285500000001 -> 670.71429
$ws.Cells.Item(98, 12).Value = 1631.0001  # L98: 1849.2 -> 1631.0001
$ws.Cells.Item(98, 13).Value = 827.28571  # M98: 825.5714499999999 -> 827.28571
$ws.Cells.Item(98, 14).Value = -4627.0001  # N98: -4845.2 -> -4627.0001
$ws.Cells.Item(114, 8).Value = 8818.6  # H114: 9687.333000000001 -> 8818.6
$ws.Cells.Item(114, 10).Value = 9773.25  # J114: 12031 -> 9773.25
$ws.Cells.Item(114, 12).Value = 29319.75  # L114: 36093 -> 29319.75
$ws.Cells.Item(114, 14).Value = -35827.75  # N114: -42601 -> -35827.75
$ws.Cells.Item(117, 8).Value = 0  # H117: 950 -> 0
$ws.Cells.Item(117, 10).Value = 0  # J117: 950 -> 0
$ws.Cells.Item(117, 12).Value = 0  # L117: 2850 -> 0
$ws.Cells.Item(117, 14).ClearContents()  # N117: -9734 -> (removed)
$ws.Cells.Item(119, 8).Value = 699  # H119: 700 -> 699
$ws.Cells.Item(119, 9).Value = 699  # I119: 700 -> 699
$ws.Cells.Item(119, 11).Value = 2097  # K119: 2100 -> 2097
$ws.Cells.Item(119, 13).Value = 2741  # M119: 2738 -> 2741
$ws.Cells.Item(121, 8).Value = 500  # H121: 531.1875 -> 500
$ws.Cells.Item(121, 10).Value = 0  # J121: 749.5 -> 0
$ws.Cells.Item(121, 12).Value = 0  # L121: 2248.5 -> 0
$ws.Cells.Item(121, 14).ClearContents()  # N121: -4868.5 -> (removed)
$ws.Cells.Item(122, 8).Value = 967.2  # H122: 1085.25 -> 967.2
$ws.Cells.Item(122, 9).Value = 622.75  # I122: 665.3333 -> 622.75
$ws.Cells.Item(122, 11).Value = 5604.75  # K122: 5987.9997 -> 5604.75
$ws.Cells.Item(122, 13).Value = -3154.75  # M122: -3537.9997 -> -3154.75
$ws.Cells.Item(124, 8).Value = 5025.5557  # H124: 5030.25 -> 5025.5557
$ws.Cells.Item(124, 9).Value = 4999.5  # I124: 5000 -> 4999.5
$ws.Cells.Item(124, 11).Value = 14998.5  # K124: 15000 -> 14998.5
$ws.Cells.Item(124, 13).Value = -10088.5  # M124: -10090 -> -10088.5
$ws.Cells.Item(125, 8).Value = 3000  # H125: 10000 -> 3000
$ws.Cells.Item(125, 9).Value = 3000  # I125: 0 -> 3000
$ws.Cells.Item(125, 10).Value = 0  # J125: 10000 -> 0
$ws.Cells.Item(125, 11).Value = 9000  # K125: 0 -> 9000
$ws.Cells.Item(125, 12).Value = 0  # L125: 30000 -> 0
$ws.Cells.Item(125, 13).Value = -4080  # M125: None -> -4080
$ws.Cells.Item(125, 14).ClearContents()  # N125: -39840 -> (removed)
$ws.Cells.Item(129, 8).Value = 1400  # H129: 1232 -> 1400
$ws.Cells.Item(129, 10).Value = 0  # J129: 1190 -> 0
$ws.Cells.Item(129, 12).Value = 0  # L129: 3570 -> 0
$ws.Cells.Item(129, 14).ClearContents()  # N129: -13570 -> (removed)
$ws.Cells.Item(131, 8).Value = 1000  # H131: 0 -> 1000
$ws.Cells.Item(131, 10).Value = 1000  # J131: 0 -> 1000
$ws.Cells.Item(131, 12).Value = 3000  # L131: 0 -> 3000
$ws.Cells.Item(131, 14).Value = -13080  # N131: None -> -13080
$ws.Cells.Item(132, 8).Value = 4739.6665  # H132: 4863.9287 -> 4739.6665
$ws.Cells.Item(132, 10).Value = 5799.8887  # J132: 6149.875 -> 5799.8887
$ws.Cells.Item(132, 12).Value = 52198.99830000001  # L132: 55348.875 -> 52198.99830000001
$ws.Cells.Item(132, 14).Value = -57258.99830000001  # N132: -60408.875 -> -57258.99830000001
$ws.Cells.Item(133, 8).Value = 4649.6665  # H133: 4414 -> 4649.6665
$ws.Cells.Item(133, 10).Value = 7499.5  # J133: 5999.6665 -> 7499.5
$ws.Cells.Item(133, 12).Value = 22498.5  # L133: 17998.9995 -> 22498.5
$ws.Cells.Item(133, 14).Value = -32618.5  # N133: -28118.9995 -> -32618.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 121.53333  # H2: 124.64286 -> 121.53333
$ws.Cells.Item(2, 10).Value = 200.28572  # J2: 220.66667 -> 200.28572
$ws.Cells.Item(2, 12).Value = 200.28572  # L2: 220.66667 -> 200.28572
$ws.Cells.Item(2, 14).Value = -426.28572  # N2: -446.66667 -> -426.28572
$ws.Cells.Item(113, 8).Value = 4332.647  # H113: 4591.3125 -> 4332.647
$ws.Cells.Item(113, 9).Value = 1165.5  # I113: 1273.4445 -> 1165.5
$ws.Cells.Item(113, 11).Value = 1165.5  # K113: 1273.4445 -> 1165.5
$ws.Cells.Item(113, 13).Value = 1004.5  # M113: 896.5554999999999 -> 1004.5
$ws.Cells.Item(132, 8).Value = 94430.14  # H132: 77707.06 -> 94430.14
$ws.Cells.Item(132, 9).Value = 100155.54  # I132: 99386.16 -> 100155.54
$ws.Cells.Item(132, 10).Value = 20000  # J132: 7250 -> 20000
$ws.Cells.Item(132, 11).Value = 300466.62  # K132: 298158.48 -> 300466.62
$ws.Cells.Item(132, 12).Value = 60000  # L132: 21750 -> 60000
$ws.Cells.Item(132, 13).Value = -297936.62  # M132: -295628.48 -> -297936.62
$ws.Cells.Item(132, 14).Value = -65060  # N132: -26810 -> -65060
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(10, 8).Value = 0  # H10: 9752.5 -> 0
$ws.Cells.Item(10, 9).Value = 0  # I10: 9752.5 -> 0
$ws.Cells.Item(10, 11).Value = 0  # K10: 9752.5 -> 0
$ws.Cells.Item(10, 13).ClearContents()  # M10: -9583.5 -> (removed)
$ws.Cells.Item(45, 8).Value = 22499  # H45: 22499.5 -> 22499
$ws.Cells.Item(45, 9).Value = 0  # I45: 14999 -> 0
$ws.Cells.Item(45, 10).Value = 22499  # J45: 30000 -> 22499
$ws.Cells.Item(45, 11).Value = 0  # K45: 14999 -> 0
$ws.Cells.Item(45, 12).Value = 22499  # L45: 30000 -> 22499
$ws.Cells.Item(45, 13).ClearContents()  # M45: -14508 -> (removed)
$ws.Cells.Item(45, 14).Value = -23481  # N45: -30982 -> -23481
$ws.Cells.Item(81, 8).Value = 4903.1665  # H81: 4903.3335 -> 4903.1665
$ws.Cells.Item(81, 10).Value = 4419  # J81: 4420 -> 4419
$ws.Cells.Item(81, 12).Value = 8838  # L81: 8840 -> 8838
$ws.Cells.Item(81, 14).Value = -10960  # N81: -10962 -> -10960
$ws.Cells.Item(84, 8).Value = 4903.1665  # H84: 4903.3335 -> 4903.1665
$ws.Cells.Item(84, 10).Value = 4419  # J84: 4420 -> 4419
$ws.Cells.Item(84, 12).Value = 44190  # L84: 44200 -> 44190
$ws.Cells.Item(84, 14).Value = -54798  # N84: -54808 -> -54798
